# Move the "Backup" slide (currently the last slide, position 16) to
# position 7, right after "Grover Oracle - Implementation" and before
# "Grover Diffuser". Then update its title text.

$p = $ppt.ActivePresentation

$backup = $p.Slides.Item(16)
$backup.MoveTo(7)

# MoveTo applies immediately, but the slide object captured before the move
# keeps addressing its construction-time position for subsequent property
# access (e.g. Shapes). Re-fetch the slide at its new position (7) to edit it.
$oracleSlide = $p.Slides.Item(7)

$titleShape = $oracleSlide.Shapes.Item(1)
$tr = $titleShape.TextFrame.TextRange
$tr.Text = "Oracle Gate " + [char]0x2013 + " "
$tr.InsertAfter("Examples")
$tr.InsertAfter(" with M > 1")
